$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Add a new "properties" worksheet at the end of the workbook (after "model")
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "properties"

# ---------------------------------------------------------------------------
# Header row (bold)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "partition"
$ws.Range("B1").Value = "aspect"
$ws.Range("C1").Value = "key"
$ws.Range("D1").Value = "type"
$ws.Range("E1").Value = "value"
$ws.Range("A1:E1").Font.Bold = $true

# ---------------------------------------------------------------------------
# Data rows
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "Table"
$ws.Range("B2").Value = "security"
$ws.Range("C2").Value = "locked"
$ws.Range("D2").Value = "boolean"
$ws.Range("E2").Formula = "=""FALSE"""

$ws.Range("A3").Value = "Table"
$ws.Range("B3").Value = "security"
$ws.Range("C3").Value = "unverifiedUserCanCreate"
$ws.Range("D3").Value = "boolean"
$ws.Range("E3").Formula = "=""TRUE"""

$ws.Range("A4").Value = "Table"
$ws.Range("B4").Value = "security"
$ws.Range("C4").Value = "defaultAccessOnCreation"
$ws.Range("D4").Value = "string"
$ws.Range("E4").Value = "HIDDEN"

# C3/C4 pick up the plain (non-themed) Calibri font used elsewhere in this
# workbook (style index 10 in the original file) - reuse it by copying the
# format from an existing cell that already carries it instead of minting a
# brand new style.
$styleSource = $wb.Worksheets.Item("calculates").Range("A2")
[void]$styleSource.Copy()
[void]$ws.Range("C3:C4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Column widths (best-fit, matches the width of the widest entry per column)
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 7.8125
$ws.Columns.Item(2).ColumnWidth = 7.1875
$ws.Columns.Item(3).ColumnWidth = 22.96875
$ws.Columns.Item(4).ColumnWidth = 7.5
$ws.Columns.Item(5).ColumnWidth = 7.03125

# ---------------------------------------------------------------------------
# Freeze the header row and set the selection like the source file
# ---------------------------------------------------------------------------
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
[void]$ws.Range("D8").Select()

# ---------------------------------------------------------------------------
# Make the new sheet the active/selected tab
# ---------------------------------------------------------------------------
[void]$ws.Activate()
